$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: goods id 13001
$ws.Range("A8").Value = 13001
$ws.Range("B8").Value = 1
$ws.Range("H8").Value = '[{"item":{"id":"10001"}}]'
$ws.Range("K8").Value = '[{"money":"800"}]'
$ws.Range("L8").Value = '[{"money":"800"}]'
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1

# Row 9: goods id 13002
$ws.Range("A9").Value = 13002
$ws.Range("B9").Value = 1
$ws.Range("H9").Value = '[{"item":{"id":"11002"}}]'
$ws.Range("K9").Value = '[{"money":"80"}]'
$ws.Range("L9").Value = '[{"money":"80"}]'
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 1

# Row 10: goods id 13003
$ws.Range("A10").Value = 13003
$ws.Range("B10").Value = 1
$ws.Range("H10").Value = '[{"item":{"id":"12001"}}]'
$ws.Range("K10").Value = '[{"money":"400"}]'
$ws.Range("L10").Value = '[{"money":"400"}]'
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 1
